$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> new values
$ws.Range("A2").Value = 112183485
$ws.Range("B2").Value = 5426
$ws.Range("E2").Value = 101410
$ws.Range("F2").Value = "Reliktbock"
$ws.Range("G2").Value = "Nothorhina muricata"
$ws.Range("H2").Value = "(Dalman, 1817)"
$ws.Range("Q2").Value = 387592
$ws.Range("R2").Value = 6855479

# Row 3 -> new values
$ws.Range("A3").Value = 112183324
$ws.Range("B3").Value = 90837
$ws.Range("E3").Value = 5966
$ws.Range("F3").Value = "Motaggsvamp"
$ws.Range("G3").Value = "Sarcodon squamosus"
$ws.Range("H3").Value = "(Schaeff.) Quél."
$ws.Range("Q3").Value = 387555
$ws.Range("R3").Value = 6855526

# Row 4 -> new values
$ws.Range("A4").Value = 112183589
$ws.Range("B4").Value = 90826
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 4366
$ws.Range("F4").Value = "Skarp dropptaggsvamp"
$ws.Range("G4").Value = "Hydnellum peckii"
$ws.Range("H4").Value = "Banker"
$ws.Range("Q4").Value = 387651
$ws.Range("R4").Value = 6855494

# Row 5 -> new values
$ws.Range("A5").Value = 112182653
$ws.Range("B5").Value = 90808
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 4362
$ws.Range("F5").Value = "Blå taggsvamp"
$ws.Range("G5").Value = "Hydnellum caeruleum"
$ws.Range("H5").Value = "(Hornem.) P.Karst."
$ws.Range("Q5").Value = 387566
$ws.Range("R5").Value = 6855527
